$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 152. Every row from 152..183 shifts down to 153..184.
$ws.Rows.Item(152).Insert()

# After the shift, the data that used to live in row 156 is now in row 157.
# The new row 152 should carry a copy of that record (same market/price/unit info)
# but with an updated sampling date.
$src = $ws.Rows.Item(157)
$dst = $ws.Rows.Item(152)
$src.Copy()
$dst.PasteSpecial(-4104)  # xlPasteAll
$excel.CutCopyMode = $false

$ws.Range("D152").Value2 = 44964
